# Add new buildings to the Improvements/Buildings list on Sheet1.
#
# Layout before the edit (column B, 1 building per row, blank rows used
# as separators between groups):
#   ...
#   20 Granary
#   21 Smokehouse
#   23 Harbor          24 Docks          25 Shipyard
#   27 Forge           28 Blacksmith     29 Fletcher
#   31 Archery Range   32 Training Yard  33 Stable   34 Siege workshop
#   36 Workshop
#
# After the edit:
#   - Two new rows (Grocer, Apotecary) are inserted right after
#     Smokehouse (row 21), pushing every following row down by 2.
#   - Two new rows (Mage guild, Mage tower) are appended after
#     Workshop, leaving a blank separator row in between, matching the
#     existing group-separation convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 23 ("Harbor"),
# shifting the whole lower block (Harbor...Workshop) down by 2.
$ws.Rows("23:24").Insert() | Out-Null

# Fill the two newly inserted rows.
$ws.Range("B22").Value = "Grocer"
$ws.Range("B23").Value = "Apotecary"

# Append the two new Mage buildings after the (now shifted) Workshop
# row (38), leaving row 39 blank as a group separator.
$ws.Range("B40").Value = "Mage guild"
$ws.Range("B41").Value = "Mage tower"

# Match the author's final selection/active cell.
$ws.Range("B42").Select() | Out-Null
